$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking columns (Price, Volume%, Hora)
# so Excel does not auto-convert these strings to numbers and lose exact formatting
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Apply cell value changes per the diff
# Row 2
$ws.Range("D2").Value = "290.72"
$ws.Range("E2").Value = "-8.07%"
$ws.Range("G2").Value = "22"

# Row 3
$ws.Range("D3").Value = "40.32"
$ws.Range("E3").Value = "-2.21%"
$ws.Range("G3").Value = "22"

# Row 4
$ws.Range("D4").Value = "5.043"
$ws.Range("E4").Value = "-3.52%"
$ws.Range("G4").Value = "22"

# Row 5
$ws.Range("D5").Value = "0.07285"
$ws.Range("E5").Value = "-4.55%"
$ws.Range("G5").Value = "22"

# Row 6
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").Value = "1.566"
$ws.Range("E6").Value = "-5.23%"
$ws.Range("G6").Value = "22"

# Row 7
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "0.9204"
$ws.Range("E7").Value = "-1.18%"
$ws.Range("G7").Value = "22"

# Row 8
$ws.Range("B8").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C8").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D8").Value = "0.1158"
$ws.Range("E8").Value = "-6.60%"
$ws.Range("G8").Value = "22"

# Row 9
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "0.1726"
$ws.Range("E9").Value = "-5.35%"
$ws.Range("G9").Value = "22"

# Row 10
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "0.08614"
$ws.Range("E10").Value = "-5.44%"
$ws.Range("G10").Value = "22"

# Row 11
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "0.04185"
$ws.Range("E11").Value = "0.63%"
$ws.Range("G11").Value = "22"

# Row 12
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").Value = "0.1053"
$ws.Range("E12").Value = "-0.09%"
$ws.Range("G12").Value = "22"

# Row 13
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").Value = "0.001268"
$ws.Range("E13").Value = "-0.24%"
$ws.Range("G13").Value = "22"

# Row 14
$ws.Range("B14").Value = "TigerCash"
$ws.Range("C14").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D14").Value = "0.005808"
$ws.Range("E14").Value = "-2.17%"
$ws.Range("G14").Value = "22"

# Row 15
$ws.Range("B15").Value = "LEO"
$ws.Range("C15").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D15").Value = "3.400"
$ws.Range("E15").Value = "1.26%"
$ws.Range("G15").Value = "22"

# Row 16
$ws.Range("B16").Value = "GateToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D16").Value = "4.284"
$ws.Range("E16").Value = "-1.31%"
$ws.Range("G16").Value = "22"

# Row 17
$ws.Range("D17").Value = "2.333"
$ws.Range("E17").Value = "-3.79%"
$ws.Range("G17").Value = "22"

# Row 18
$ws.Range("E18").Value = "-2.50%"
$ws.Range("G18").Value = "22"

# Row 19
$ws.Range("D19").Value = "7.850"
$ws.Range("E19").Value = "-7.19%"
$ws.Range("G19").Value = "22"

# Row 20
$ws.Range("D20").Value = "0.1381"
$ws.Range("E20").Value = "1.51%"
$ws.Range("G20").Value = "22"

# Row 21
$ws.Range("D21").Value = "0.2885"
$ws.Range("G21").Value = "22"

# Row 22
$ws.Range("D22").Value = "0.03860"
$ws.Range("E22").Value = "-4.51%"
$ws.Range("G22").Value = "22"

# Row 23
$ws.Range("D23").Value = "0.001271"
$ws.Range("E23").Value = "-0.43%"
$ws.Range("G23").Value = "22"

# Row 24
$ws.Range("D24").Value = "0.003813"
$ws.Range("E24").Value = "-6.63%"
$ws.Range("G24").Value = "22"

# Row 25
$ws.Range("D25").Value = "0.0001282"
$ws.Range("E25").Value = "0.11%"
$ws.Range("G25").Value = "22"

# Row 26
$ws.Range("D26").Value = "0.0003726"
$ws.Range("G26").Value = "22"

# Row 27
$ws.Range("G27").Value = "22"

# Row 28
$ws.Range("G28").Value = "22"

# Row 29
$ws.Range("G29").Value = "22"

# Row 30
$ws.Range("G30").Value = "22"

# Row 31
$ws.Range("G31").Value = "22"

# Row 32
$ws.Range("G32").Value = "22"

# Row 33
$ws.Range("G33").Value = "22"

# Row 34
$ws.Range("G34").Value = "22"

# Row 35
$ws.Range("G35").Value = "22"

# Row 36
$ws.Range("G36").Value = "22"

# Row 37
$ws.Range("G37").Value = "22"

# Row 38
$ws.Range("E38").Value = "-5.64%"
$ws.Range("G38").Value = "22"

# Row 39
$ws.Range("D39").Value = "0.04951"
$ws.Range("E39").Value = "-5.26%"
$ws.Range("G39").Value = "22"

# Row 40
$ws.Range("D40").Value = "0.006638"
$ws.Range("E40").Value = "206.27%"
$ws.Range("G40").Value = "22"

# Row 41
$ws.Range("D41").Value = "0.007678"
$ws.Range("E41").Value = "-1.55%"
$ws.Range("G41").Value = "22"

# Row 42
$ws.Range("D42").Value = "0.1272"
$ws.Range("E42").Value = "-1.72%"
$ws.Range("G42").Value = "22"

# Row 43
$ws.Range("D43").Value = "0.007371"
$ws.Range("E43").Value = "3.93%"
$ws.Range("G43").Value = "22"

# Row 44
$ws.Range("D44").Value = "0.007066"
$ws.Range("E44").Value = "-14.26%"
$ws.Range("G44").Value = "22"

# Row 45
$ws.Range("D45").Value = "0.2905"
$ws.Range("E45").Value = "-15.48%"
$ws.Range("G45").Value = "22"

# Row 46
$ws.Range("D46").Value = "0.00006417"
$ws.Range("E46").Value = "-4.22%"
$ws.Range("G46").Value = "22"

# Row 47
$ws.Range("E47").Value = "-0.60%"
$ws.Range("G47").Value = "22"

# Row 48
$ws.Range("D48").Value = "0.01737"
$ws.Range("E48").Value = "-95.08%"
$ws.Range("G48").Value = "22"

# Row 49
$ws.Range("E49").Value = "-0.68%"
$ws.Range("G49").Value = "22"

# Row 50
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").Value = "-0.60%"
$ws.Range("G50").Value = "22"

# Row 51
$ws.Range("E51").Value = "-0.60%"
$ws.Range("G51").Value = "22"

